# Update NATMI LR-pair metrics for Lgals3bp-Itgb1 (Young D4) per Dr Hou advice.
# Ligand/receptor-expressing cell counts go from 1 -> 3, with all dependent
# expression/specificity metrics recomputed accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 77.88937866666667
$ws.Range("H2").Value = 233.668136
$ws.Range("I2").Value = 0.117745379143944
$ws.Range("J2").Value = 0.117745379143944
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 163.7119853333333
$ws.Range("N2").Value = 491.135956
$ws.Range("O2").Value = 0.2754003062401033
$ws.Range("P2").Value = 0.2754003062401033
$ws.Range("Q2").Value = 12751.42481789978
$ws.Range("R2").Value = 114762.823361098
$ws.Range("S2").Value = 0.03242711347459924
$ws.Range("T2").Value = 0.03242711347459924

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 77.88937866666667
$ws.Range("H3").Value = 233.668136
$ws.Range("I3").Value = 0.117745379143944
$ws.Range("J3").Value = 0.117745379143944
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 172.558497
$ws.Range("N3").Value = 517.675491
$ws.Range("O3").Value = 0.290282124557779
$ws.Range("P3").Value = 0.290282124557779
$ws.Range("Q3").Value = 13440.47411498386
$ws.Range("R3").Value = 120964.2670348548
$ws.Range("S3").Value = 0.03417937881476526
$ws.Range("T3").Value = 0.03417937881476526

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 77.88937866666667
$ws.Range("H4").Value = 233.668136
$ws.Range("I4").Value = 0.117745379143944
$ws.Range("J4").Value = 0.117745379143944
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 111.4881643333333
$ws.Range("N4").Value = 334.464493
$ws.Range("O4").Value = 0.1875481171218523
$ws.Range("P4").Value = 0.1875481171218523
$ws.Range("Q4").Value = 8683.743848610562
$ws.Range("R4").Value = 78153.69463749505
$ws.Range("S4").Value = 0.02208292415824531
$ws.Range("T4").Value = 0.02208292415824531

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 77.88937866666667
$ws.Range("H5").Value = 233.668136
$ws.Range("I5").Value = 0.117745379143944
$ws.Range("J5").Value = 0.117745379143944
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 86.95798233333333
$ws.Range("N5").Value = 260.873947
$ws.Range("O5").Value = 0.1462828449356383
$ws.Range("P5").Value = 0.1462828449356383
$ws.Range("Q5").Value = 6773.103214050311
$ws.Range("R5").Value = 60957.92892645279
$ws.Range("S5").Value = 0.0172241290392015
$ws.Range("T5").Value = 0.0172241290392015

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 77.88937866666667
$ws.Range("H6").Value = 233.668136
$ws.Range("I6").Value = 0.117745379143944
$ws.Range("J6").Value = 0.117745379143944
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 59.73436333333333
$ws.Range("N6").Value = 179.20309
$ws.Range("O6").Value = 0.100486607144627
$ws.Range("P6").Value = 0.100486607144627
$ws.Range("Q6").Value = 4652.672445082248
$ws.Range("R6").Value = 41874.05200574023
$ws.Range("S6").Value = 0.01183183365713265
$ws.Range("T6").Value = 0.01183183365713265

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 221.6758116666666
$ws.Range("H7").Value = 665.027435
$ws.Range("I7").Value = 0.3351073399036296
$ws.Range("J7").Value = 0.3351073399036297
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 163.7119853333333
$ws.Range("N7").Value = 491.135956
$ws.Range("O7").Value = 0.2754003062401033
$ws.Range("P7").Value = 0.2754003062401033
$ws.Range("Q7").Value = 36290.98722832809
$ws.Range("R7").Value = 326618.8850549528
$ws.Range("S7").Value = 0.09228866403276599
$ws.Range("T7").Value = 0.092288664032766

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 221.6758116666666
$ws.Range("H8").Value = 665.027435
$ws.Range("I8").Value = 0.3351073399036296
$ws.Range("J8").Value = 0.3351073399036297
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 172.558497
$ws.Range("N8").Value = 517.675491
$ws.Range("O8").Value = 0.290282124557779
$ws.Range("P8").Value = 0.290282124557779
$ws.Range("Q8").Value = 38252.04488245506
$ws.Range("R8").Value = 344268.4039420955
$ws.Range("S8").Value = 0.09727567058213138
$ws.Range("T8").Value = 0.09727567058213141

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 221.6758116666666
$ws.Range("H9").Value = 665.027435
$ws.Range("I9").Value = 0.3351073399036296
$ws.Range("J9").Value = 0.3351073399036297
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 111.4881643333333
$ws.Range("N9").Value = 334.464493
$ws.Range("O9").Value = 0.1875481171218523
$ws.Range("P9").Value = 0.1875481171218523
$ws.Range("Q9").Value = 24714.22931981838
$ws.Range("R9").Value = 222428.0638783654
$ws.Range("S9").Value = 0.0628487506326383
$ws.Range("T9").Value = 0.06284875063263831

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 221.6758116666666
$ws.Range("H10").Value = 665.027435
$ws.Range("I10").Value = 0.3351073399036296
$ws.Range("J10").Value = 0.3351073399036297
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 86.95798233333333
$ws.Range("N10").Value = 260.873947
$ws.Range("O10").Value = 0.1462828449356383
$ws.Range("P10").Value = 0.1462828449356383
$ws.Range("Q10").Value = 19276.48131463733
$ws.Range("R10").Value = 173488.3318317359
$ws.Range("S10").Value = 0.0490204550399169
$ws.Range("T10").Value = 0.04902045503991691

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 221.6758116666666
$ws.Range("H11").Value = 665.027435
$ws.Range("I11").Value = 0.3351073399036296
$ws.Range("J11").Value = 0.3351073399036297
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 59.73436333333333
$ws.Range("N11").Value = 179.20309
$ws.Range("O11").Value = 0.100486607144627
$ws.Range("P11").Value = 0.100486607144627
$ws.Range("Q11").Value = 13241.66347630824
$ws.Range("R11").Value = 119174.9712867741
$ws.Range("S11").Value = 0.03367379961617701
$ws.Range("T11").Value = 0.03367379961617701

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 131.496086
$ws.Range("H12").Value = 394.488258
$ws.Range("I12").Value = 0.1987826423455699
$ws.Range("J12").Value = 0.1987826423455699
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 163.7119853333333
$ws.Range("N12").Value = 491.135956
$ws.Range("O12").Value = 0.2754003062401033
$ws.Range("P12").Value = 0.2754003062401033
$ws.Range("Q12").Value = 21527.48530262274
$ws.Range("R12").Value = 193747.3677236046
$ws.Range("S12").Value = 0.05474480057718688
$ws.Range("T12").Value = 0.05474480057718688

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 131.496086
$ws.Range("H13").Value = 394.488258
$ws.Range("I13").Value = 0.1987826423455699
$ws.Range("J13").Value = 0.1987826423455699
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 172.558497
$ws.Range("N13").Value = 517.675491
$ws.Range("O13").Value = 0.290282124557779
$ws.Range("P13").Value = 0.290282124557779
$ws.Range("Q13").Value = 22690.76696154274
$ws.Range("R13").Value = 204216.9026538846
$ws.Range("S13").Value = 0.05770304774528114
$ws.Range("T13").Value = 0.05770304774528116

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 131.496086
$ws.Range("H14").Value = 394.488258
$ws.Range("I14").Value = 0.1987826423455699
$ws.Range("J14").Value = 0.1987826423455699
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 111.4881643333333
$ws.Range("N14").Value = 334.464493
$ws.Range("O14").Value = 0.1875481171218523
$ws.Range("P14").Value = 0.1875481171218523
$ws.Range("Q14").Value = 14660.25724515813
$ws.Range("R14").Value = 131942.3152064232
$ws.Range("S14").Value = 0.03728131028841822
$ws.Range("T14").Value = 0.03728131028841822

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 131.496086
$ws.Range("H15").Value = 394.488258
$ws.Range("I15").Value = 0.1987826423455699
$ws.Range("J15").Value = 0.1987826423455699
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 86.95798233333333
$ws.Range("N15").Value = 260.873947
$ws.Range("O15").Value = 0.1462828449356383
$ws.Range("P15").Value = 0.1462828449356383
$ws.Range("Q15").Value = 11434.63432329048
$ws.Range("R15").Value = 102911.7089096143
$ws.Range("S15").Value = 0.02907849044613346
$ws.Range("T15").Value = 0.02907849044613346

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 131.496086
$ws.Range("H16").Value = 394.488258
$ws.Range("I16").Value = 0.1987826423455699
$ws.Range("J16").Value = 0.1987826423455699
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 59.73436333333333
$ws.Range("N16").Value = 179.20309
$ws.Range("O16").Value = 0.100486607144627
$ws.Range("P16").Value = 0.100486607144627
$ws.Range("Q16").Value = 7854.834978035246
$ws.Range("R16").Value = 70693.51480231721
$ws.Range("S16").Value = 0.01997499328855017
$ws.Range("T16").Value = 0.01997499328855017

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 224.0923716666667
$ws.Range("H17").Value = 672.277115
$ws.Range("I17").Value = 0.3387604538235877
$ws.Range("J17").Value = 0.3387604538235878
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 163.7119853333333
$ws.Range("N17").Value = 491.135956
$ws.Range("O17").Value = 0.2754003062401033
$ws.Range("P17").Value = 0.2754003062401033
$ws.Range("Q17").Value = 36686.60706360521
$ws.Range("R17").Value = 330179.4635724469
$ws.Range("S17").Value = 0.09329473272505243
$ws.Range("T17").Value = 0.09329473272505245

# Row 18
$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 224.0923716666667
$ws.Range("H18").Value = 672.277115
$ws.Range("I18").Value = 0.3387604538235877
$ws.Range("J18").Value = 0.3387604538235878
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 172.558497
$ws.Range("N18").Value = 517.675491
$ws.Range("O18").Value = 0.290282124557779
$ws.Range("P18").Value = 0.290282124557779
$ws.Range("Q18").Value = 38669.04284396538
$ws.Range("R18").Value = 348021.3855956884
$ws.Range("S18").Value = 0.09833610425206842
$ws.Range("T18").Value = 0.09833610425206843

# Row 19
$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 224.0923716666667
$ws.Range("H19").Value = 672.277115
$ws.Range("I19").Value = 0.3387604538235877
$ws.Range("J19").Value = 0.3387604538235878
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 111.4881643333333
$ws.Range("N19").Value = 334.464493
$ws.Range("O19").Value = 0.1875481171218523
$ws.Range("P19").Value = 0.1875481171218523
$ws.Range("Q19").Value = 24983.64715821974
$ws.Range("R19").Value = 224852.8244239777
$ws.Range("S19").Value = 0.06353388526995807
$ws.Range("T19").Value = 0.06353388526995808

# Row 20
$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 224.0923716666667
$ws.Range("H20").Value = 672.277115
$ws.Range("I20").Value = 0.3387604538235877
$ws.Range("J20").Value = 0.3387604538235878
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 86.95798233333333
$ws.Range("N20").Value = 260.873947
$ws.Range("O20").Value = 0.1462828449356383
$ws.Range("P20").Value = 0.1462828449356383
$ws.Range("Q20").Value = 19486.62049642477
$ws.Range("R20").Value = 175379.5844678229
$ws.Range("S20").Value = 0.04955484293700235
$ws.Range("T20").Value = 0.04955484293700236

# Row 21
$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 224.0923716666667
$ws.Range("H21").Value = 672.277115
$ws.Range("I21").Value = 0.3387604538235877
$ws.Range("J21").Value = 0.3387604538235878
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 59.73436333333333
$ws.Range("N21").Value = 179.20309
$ws.Range("O21").Value = 0.100486607144627
$ws.Range("P21").Value = 0.100486607144627
$ws.Range("Q21").Value = 13386.01514936504
$ws.Range("R21").Value = 120474.1363442853
$ws.Range("S21").Value = 0.0340408886395064
$ws.Range("T21").Value = 0.03404088863950641

# Row 22
$ws.Range("E22").Value = 3
$ws.Range("G22").Value = 6.353234333333333
$ws.Range("H22").Value = 19.059703
$ws.Range("I22").Value = 0.009604184783268722
$ws.Range("J22").Value = 0.009604184783268722
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 163.7119853333333
$ws.Range("N22").Value = 491.135956
$ws.Range("O22").Value = 0.2754003062401033
$ws.Range("P22").Value = 0.2754003062401033
$ws.Range("Q22").Value = 1040.100605997896
$ws.Range("R22").Value = 9360.905453981068
$ws.Range("S22").Value = 0.002644995430498746
$ws.Range("T22").Value = 0.002644995430498746

# Row 23
$ws.Range("E23").Value = 3
$ws.Range("G23").Value = 6.353234333333333
$ws.Range("H23").Value = 19.059703
$ws.Range("I23").Value = 0.009604184783268722
$ws.Range("J23").Value = 0.009604184783268722
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 172.558497
$ws.Range("N23").Value = 517.675491
$ws.Range("O23").Value = 0.290282124557779
$ws.Range("P23").Value = 0.290282124557779
$ws.Range("Q23").Value = 1096.304567648797
$ws.Range("R23").Value = 9866.741108839171
$ws.Range("S23").Value = 0.002787923163532737
$ws.Range("T23").Value = 0.002787923163532737

# Row 24
$ws.Range("E24").Value = 3
$ws.Range("G24").Value = 6.353234333333333
$ws.Range("H24").Value = 19.059703
$ws.Range("I24").Value = 0.009604184783268722
$ws.Range("J24").Value = 0.009604184783268722
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 111.4881643333333
$ws.Range("N24").Value = 334.464493
$ws.Range("O24").Value = 0.1875481171218523
$ws.Range("P24").Value = 0.1875481171218523
$ws.Range("Q24").Value = 708.310433402842
$ws.Range("R24").Value = 6374.793900625578
$ws.Range("S24").Value = 0.001801246772592394
$ws.Range("T24").Value = 0.001801246772592394

# Row 25
$ws.Range("E25").Value = 3
$ws.Range("G25").Value = 6.353234333333333
$ws.Range("H25").Value = 19.059703
$ws.Range("I25").Value = 0.009604184783268722
$ws.Range("J25").Value = 0.009604184783268722
$ws.Range("K25").Value = 3
$ws.Range("M25").Value = 86.95798233333333
$ws.Range("N25").Value = 260.873947
$ws.Range("O25").Value = 0.1462828449356383
$ws.Range("P25").Value = 0.1462828449356383
$ws.Range("Q25").Value = 552.4644389175268
$ws.Range("R25").Value = 4972.179950257741
$ws.Range("S25").Value = 0.001404927473384116
$ws.Range("T25").Value = 0.001404927473384116

# Row 26
$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 6.353234333333333
$ws.Range("H26").Value = 19.059703
$ws.Range("I26").Value = 0.009604184783268722
$ws.Range("J26").Value = 0.009604184783268722
$ws.Range("K26").Value = 3
$ws.Range("M26").Value = 59.73436333333333
$ws.Range("N26").Value = 179.20309
$ws.Range("O26").Value = 0.100486607144627
$ws.Range("P26").Value = 0.100486607144627
$ws.Range("Q26").Value = 379.5064080091411
$ws.Range("R26").Value = 3415.557672082269
$ws.Range("S26").Value = 0.0009650919432607284
$ws.Range("T26").Value = 0.0009650919432607284
